$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while forcing text storage
# (so numeric-looking strings like "1.00" are not coerced to numbers),
# then restore the cell style so no stray formatting is introduced.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '64.345.11'
Set-TextValue $ws.Range("E2") '  -0.37%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.140.05'
Set-TextValue $ws.Range("E3") '  -0.44%  '

# Row 4
Set-TextValue $ws.Range("E4") '  +0.03%  '

# Row 5
Set-TextValue $ws.Range("D5") '613.49'
Set-TextValue $ws.Range("E5") '  +1.15%  '

# Row 6
Set-TextValue $ws.Range("D6") '143.06'
Set-TextValue $ws.Range("E6") '  -3.89%  '

# Row 7
Set-TextValue $ws.Range("E7") '  -0.02%  '

# Row 8
Set-TextValue $ws.Range("D8") '3.142.31'
Set-TextValue $ws.Range("E8") '  -0.26%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.521'
Set-TextValue $ws.Range("E9") '  -1.53%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.150'
Set-TextValue $ws.Range("E10") '  -1.92%  '

# Row 11
Set-TextValue $ws.Range("D11") '5.37'
Set-TextValue $ws.Range("E11") '  -3.88%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.466'
Set-TextValue $ws.Range("E12") '  -2.50%  '

# Row 13
Set-TextValue $ws.Range("D13") '0.0000254'
Set-TextValue $ws.Range("E13") '  -2.05%  '

# Row 14
Set-TextValue $ws.Range("D14") '35.23'
Set-TextValue $ws.Range("E14") '  -3.72%  '

# Row 15
Set-TextValue $ws.Range("D15") '3.656.13'
Set-TextValue $ws.Range("E15") '  -0.48%  '

# Row 16
Set-TextValue $ws.Range("E16") '  +2.98%  '

# Row 17
Set-TextValue $ws.Range("D17") '64.259.94'
Set-TextValue $ws.Range("E17") '  -0.49%  '

# Row 18
Set-TextValue $ws.Range("D18") '3.141.73'
Set-TextValue $ws.Range("E18") '  -0.47%  '

# Row 19
Set-TextValue $ws.Range("D19") '6.80'
Set-TextValue $ws.Range("E19") '  -2.45%  '

# Row 20
Set-TextValue $ws.Range("D20") '475.25'
Set-TextValue $ws.Range("E20") '  -1.75%  '

# Row 21
Set-TextValue $ws.Range("D21") '14.52'
Set-TextValue $ws.Range("E21") '  -0.57%  '

# Row 22
Set-TextValue $ws.Range("D22") '0.717'
Set-TextValue $ws.Range("E22") '  +0.83%  '

# Row 23
Set-TextValue $ws.Range("D23") '7.84'
Set-TextValue $ws.Range("E23") '  +0.85%  '

# Row 24
Set-TextValue $ws.Range("D24") '13.64'
Set-TextValue $ws.Range("E24") '  -1.17%  '

# Row 25
Set-TextValue $ws.Range("D25") '84.44'
Set-TextValue $ws.Range("E25") '  +0.99%  '

# Row 26
Set-TextValue $ws.Range("D26") '0.999'
Set-TextValue $ws.Range("E26") '  -0.02%  '

# Row 27
Set-TextValue $ws.Range("D27") '2.79'
Set-TextValue $ws.Range("E27") '  -4.25%  '

# Row 28
Set-TextValue $ws.Range("D28") '8.46'
Set-TextValue $ws.Range("E28") '  -0.73%  '

# Row 29
Set-TextValue $ws.Range("B29") 'Hedera'
Set-TextValue $ws.Range("C29") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D29") '0.122'
Set-TextValue $ws.Range("E29") '  -2.90%  '

# Row 30
Set-TextValue $ws.Range("B30") 'NEARProtocol'
Set-TextValue $ws.Range("C30") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D30") '7.14'
Set-TextValue $ws.Range("E30") '  +3.40%  '

# Row 31
Set-TextValue $ws.Range("D31") '2.08'
Set-TextValue $ws.Range("E31") '  -6.51%  '

# Row 32
Set-TextValue $ws.Range("E32") '  +0.03%  '

# Row 33
Set-TextValue $ws.Range("D33") '26.31'
Set-TextValue $ws.Range("E33") '  -0.28%  '

# Row 34
Set-TextValue $ws.Range("D34") '2.60'
Set-TextValue $ws.Range("E34") '  -5.39%  '

# Row 35
Set-TextValue $ws.Range("D35") '1.11'
Set-TextValue $ws.Range("E35") '  +0.40%  '

# Row 36
Set-TextValue $ws.Range("B36") 'Filecoin'
Set-TextValue $ws.Range("C36") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D36") '5.92'
Set-TextValue $ws.Range("E36") '  -2.70%  '

# Row 37
Set-TextValue $ws.Range("B37") 'PEPE'
Set-TextValue $ws.Range("C37") 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D37") '0.0₃0761'
Set-TextValue $ws.Range("E37") '  +2.75%  '

# Row 38
Set-TextValue $ws.Range("D38") '52.83'
Set-TextValue $ws.Range("E38") '  -3.07%  '

# Row 39
Set-TextValue $ws.Range("D39") '3.07'
Set-TextValue $ws.Range("E39") '  -0.59%  '

# Row 40
Set-TextValue $ws.Range("D40") '452.40'
Set-TextValue $ws.Range("E40") '  -1.42%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.0392'
Set-TextValue $ws.Range("E41") '  -2.14%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.118'
Set-TextValue $ws.Range("E42") '  -4.76%  '

# Row 43
Set-TextValue $ws.Range("D43") '8.26'
Set-TextValue $ws.Range("E43") '  -2.39%  '

# Row 44
Set-TextValue $ws.Range("D44") '2.819.92'
Set-TextValue $ws.Range("E44") '  -2.40%  '

# Row 45
Set-TextValue $ws.Range("D45") '2.27'
Set-TextValue $ws.Range("E45") '  -0.60%  '

# Row 46
Set-TextValue $ws.Range("D46") '0.263'
Set-TextValue $ws.Range("E46") '  -3.12%  '

# Row 47
Set-TextValue $ws.Range("D47") '2.43'
Set-TextValue $ws.Range("E47") '  +4.27%  '

# Row 48
Set-TextValue $ws.Range("B48") 'USDe'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D48") '1.00'
Set-TextValue $ws.Range("E48") '  +0.11%  '

# Row 49
Set-TextValue $ws.Range("B49") 'InjectiveProtocol'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D49") '26.28'
Set-TextValue $ws.Range("E49") '  -1.22%  '

# Row 50
Set-TextValue $ws.Range("B50") 'Stellar'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D50") '0.113'
Set-TextValue $ws.Range("E50") '  -1.52%  '

# Row 51
Set-TextValue $ws.Range("B51") 'Arweave'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue $ws.Range("D51") '34.65'
Set-TextValue $ws.Range("E51") '  +4.62%  '
